{"js": "const styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\nlet normal = styles.items[0];\nnormal.font.color = null;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\ntry {\n  $d.Styles(\"Normal\").Delete()\n  Write-Output \"deleted\"\n} catch {\n  Write-Output \"ERR $_\"\n}\n"}
